# Add a new event row (row 22) to the "Card19" worksheet, matching the
# existing table layout (header in row 1, data rows 2.. with columns
# A..O = card, Min_Tones, Max_Tones, Tones, Revolving flats(x),
# 1.carding elemnt(o), licker_in carding element(o), Doffer carding
# element(o), cylinder(X), doffer(X), Revolving flats(o), Date, Event,
# Correction, Serviced by).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

$newRow = 22

# Column A: card number, kept as text like the rest of column A ("19")
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "19"
$ws.Cells.Item($newRow, 1).Style = "Normal"

# Columns B..K: left blank for this event (no measurement data recorded)
for ($col = 2; $col -le 11; $col++) {
    $ws.Cells.Item($newRow, $col).NumberFormat = "@"
    $ws.Cells.Item($newRow, $col).Value = ""
    $ws.Cells.Item($newRow, $col).Style = "Normal"
}

# Column L: Date
$ws.Cells.Item($newRow, 12).Value = "17/12/2025"
# Column M: Event
$ws.Cells.Item($newRow, 13).Value = "سيرفيس"
# Column N: Correction
$ws.Cells.Item($newRow, 14).Value = "تم تغير الفلاتس المتحركة وتغير اول جريده 240"
# Column O: Serviced by
$ws.Cells.Item($newRow, 15).Value = "م محمد عبدالله ،تيم الكرد"
